$d = $word.ActiveDocument

$replacements = @(
    @("2024-12-06 Friday", "2024-12-07 Saturday"),
    @("771×6=4626", "673×9=6057"),
    @("465×6=2790", "976×2=1952"),
    @("715×7=5005", "379×5=1895"),
    @("408×7=2856", "909×9=8181"),
    @("932×8=7456", "565×3=1695"),
    @("744×2=1488", "758×2=1516"),
    @("622×9=5598", "813×6=4878"),
    @("701×7=4907", "714×6=4284"),
    @("156×3=468", "557×3=1671"),
    @("693×8=5544", "357×4=1428"),
    @("518×2=1036", "214×6=1284"),
    @("424×9=3816", "502×6=3012"),
    @("501×7=3507", "113×4=452"),
    @("700×2=1400", "288×4=1152"),
    @("763×4=3052", "793×5=3965"),
    @("754×2=1508", "138×2=276"),
    @("807×7=5649", "275×9=2475"),
    @("317×9=2853", "448×7=3136"),
    @("112×2=224", "404×8=3232"),
    @("421×9=3789", "267×4=1068"),
    @("921×6=5526", "125×9=1125"),
    @("514×3=1542", "769×8=6152"),
    @("591×2=1182", "853×8=6824"),
    @("322×5=1610", "913×8=7304"),
    @("914×8=7312", "575×3=1725")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
